# Applies numeric-value updates to the Leve profit tables across all job sheets.
# Generated from the authoritative cell-level diff of Famfrit_Profits.xlsx.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L40").Value = 2714.2856
$ws.Range("H40").Value = 2619.7
$ws.Range("J40").Value = 2714.2856
$ws.Range("N40").Value = -3064.2856
$ws.Range("K98").Value = 499.5
$ws.Range("N98").Value = -3996
$ws.Range("H98").Value = 666.3333
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 998.5
$ws.Range("I98").Value = 499.5
$ws.Range("J98").Value = 1000
$ws.Range("L122").Value = 3000
$ws.Range("H122").Value = 666.3333
$ws.Range("I122").Value = 499.5
$ws.Range("J122").Value = 1000
$ws.Range("N122").Value = -7900
$ws.Range("K122").Value = 1498.5
$ws.Range("M122").Value = 951.5
$ws.Range("K132").Value = 8544.231
$ws.Range("I132").Value = 2848.077
$ws.Range("L132").Value = 1799.0001
$ws.Range("M132").Value = -6014.231
$ws.Range("J132").Value = 599.6667
$ws.Range("N132").Value = -6859.0001
$ws.Range("H132").Value = 2687.476
$ws.Range("M135").Value = -2694.8568
$ws.Range("J135").Value = 1000
$ws.Range("I135").Value = 581.0952
$ws.Range("H135").Value = 617.5217
$ws.Range("N135").Value = -14070
$ws.Range("L135").Value = 9000
$ws.Range("K135").Value = 5229.8568
$ws.Range("K138").Value = 3147.3333
$ws.Range("I138").Value = 1049.1111
$ws.Range("H138").Value = 29416854
$ws.Range("M138").Value = 1992.6667

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K2").Value = 2452.8235
$ws.Range("I2").Value = 2452.8235
$ws.Range("M2").Value = -2339.8235
$ws.Range("H2").Value = 3095.238
$ws.Range("H32").Value = 5069.825
$ws.Range("M32").Value = -4910.1353
$ws.Range("K32").Value = 5197.1353
$ws.Range("I32").Value = 5197.1353
$ws.Range("L45").Value = 6583
$ws.Range("H45").Value = 4821.2856
$ws.Range("N45").Value = -7337
$ws.Range("J45").Value = 6583
$ws.Range("K61").Value = 2578.9167
$ws.Range("H61").Value = 2539.0625
$ws.Range("I61").Value = 2578.9167
$ws.Range("M61").Value = -2366.9167
$ws.Range("M74").Value = -171003
$ws.Range("H74").Value = 171877
$ws.Range("I74").Value = 171877
$ws.Range("K74").Value = 171877
$ws.Range("H77").Value = 171877
$ws.Range("K77").Value = 859385
$ws.Range("M77").Value = -855017
$ws.Range("I77").Value = 171877
$ws.Range("I110").Value = 29909.584
$ws.Range("H110").Value = 29909.584
$ws.Range("K110").Value = 29909.584
$ws.Range("M110").Value = -27864.584
$ws.Range("M116").Value = -158.8235
$ws.Range("I116").Value = 2452.8235
$ws.Range("K116").Value = 2452.8235
$ws.Range("H116").Value = 3095.238
$ws.Range("L122").Value = 12812.571
$ws.Range("H122").Value = 3359.2942
$ws.Range("I122").Value = 2721.2
$ws.Range("J122").Value = 4270.857
$ws.Range("N122").Value = -17712.571
$ws.Range("K122").Value = 8163.599999999999
$ws.Range("M122").Value = -5713.599999999999
$ws.Range("K132").Value = 42188.346
$ws.Range("I132").Value = 14062.782
$ws.Range("M132").Value = -39658.346
$ws.Range("H132").Value = 69117.75999999999
$ws.Range("M136").Value = -5186.750100000001
$ws.Range("K136").Value = 7736.750100000001
$ws.Range("I136").Value = 2578.9167
$ws.Range("H136").Value = 2539.0625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3095.238
$ws.Range("K3").Value = 2452.8235
$ws.Range("I3").Value = 2452.8235
$ws.Range("M3").Value = -2338.8235
$ws.Range("I105").Value = 15327.467
$ws.Range("K105").Value = 15327.467
$ws.Range("H105").Value = 8209.148999999999
$ws.Range("M105").Value = -13580.467
$ws.Range("I107").Value = 1454.8334
$ws.Range("H107").Value = 1563.8536
$ws.Range("M107").Value = 465.1666
$ws.Range("K107").Value = 1454.8334
$ws.Range("M134").Value = -282
$ws.Range("K134").Value = 2817
$ws.Range("I134").Value = 939
$ws.Range("H134").Value = 3140.818

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 2142
$ws.Range("H16").Value = 2285
$ws.Range("M16").Value = -1855
$ws.Range("I16").Value = 2142
$ws.Range("I113").Value = 2142
$ws.Range("M113").Value = 28
$ws.Range("K113").Value = 2142
$ws.Range("H113").Value = 2285
$ws.Range("J131").Value = 50133
$ws.Range("N131").Value = -60213
$ws.Range("L131").Value = 50133
$ws.Range("H131").Value = 50133
$ws.Range("K132").Value = 170499.81
$ws.Range("I132").Value = 56833.27
$ws.Range("L132").Value = 12150
$ws.Range("M132").Value = -167969.81
$ws.Range("J132").Value = 4050
$ws.Range("N132").Value = -17210
$ws.Range("H132").Value = 54126.438
$ws.Range("M134").Value = -2022.75
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 4557.75
$ws.Range("L134").Value = 11250
$ws.Range("N134").Value = -16320
$ws.Range("I134").Value = 1519.25
$ws.Range("H134").Value = 1965.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M121").Value = 275
$ws.Range("K121").Value = 1035
$ws.Range("H121").Value = 345
$ws.Range("I121").Value = 345
$ws.Range("L122").Value = 46796.4
$ws.Range("H122").Value = 4999.3335
$ws.Range("I122").Value = 3998
$ws.Range("J122").Value = 5199.6
$ws.Range("N122").Value = -51696.4
$ws.Range("K122").Value = 35982
$ws.Range("M122").Value = -33532
$ws.Range("M128").Value = -474705
$ws.Range("K128").Value = 479685
$ws.Range("H128").Value = 159895
$ws.Range("I128").Value = 159895
$ws.Range("J131").Value = 1718.1333
$ws.Range("N131").Value = -15234.3999
$ws.Range("L131").Value = 5154.3999
$ws.Range("H131").Value = 1434.3793

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L102").Value = 5730.7856
$ws.Range("I102").Value = 1861.5
$ws.Range("J102").Value = 5730.7856
$ws.Range("N102").Value = -8974.785599999999
$ws.Range("K102").Value = 1861.5
$ws.Range("M102").Value = -239.5
$ws.Range("H102").Value = 3366.2222
$ws.Range("H126").Value = 4108.6665
$ws.Range("I126").Value = 3666.5
$ws.Range("M126").Value = -8529.5
$ws.Range("K126").Value = 10999.5
$ws.Range("K132").Value = 4477.5882
$ws.Range("I132").Value = 1492.5294
$ws.Range("M132").Value = -1947.5882
$ws.Range("H132").Value = 1499.1052

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L16").Value = 2868.6667
$ws.Range("J16").Value = 2868.6667
$ws.Range("N16").Value = -3208.6667
$ws.Range("H16").Value = 1623.8485
$ws.Range("L40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("N46").Value = -4775.8
$ws.Range("J46").Value = 4399.8
$ws.Range("L46").Value = 4399.8
$ws.Range("H46").Value = 2618.9524
$ws.Range("K61").Value = 2011.7142
$ws.Range("H61").Value = 2120.889
$ws.Range("I61").Value = 2011.7142
$ws.Range("M61").Value = -1809.7142
$ws.Range("I113").Value = 2011.7142
$ws.Range("M113").Value = 158.2858000000001
$ws.Range("K113").Value = 2011.7142
$ws.Range("H113").Value = 2120.889
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I57").Value = 0
$ws.Range("N57").Value = -69882.5
$ws.Range("H57").Value = 68374.5
$ws.Range("J57").Value = 68374.5
$ws.Range("L57").Value = 68374.5
$ws.Range("K57").Value = 0
$ws.Range("L122").Value = 29250
$ws.Range("H122").Value = 43803.4
$ws.Range("I122").Value = 50289.76
$ws.Range("J122").Value = 9750
$ws.Range("N122").Value = -34150
$ws.Range("K122").Value = 150869.28
$ws.Range("M122").Value = -148419.28
$ws.Range("J126").Value = 3175
$ws.Range("N126").Value = -14465
$ws.Range("H126").Value = 18526964
$ws.Range("I126").Value = 33345994
$ws.Range("M126").Value = -100035512
$ws.Range("K126").Value = 100037982
$ws.Range("L126").Value = 9525
$ws.Range("M57").ClearContents()
